# Add a new "WTI-spot" row (row 15) to the "Simplified contract multiplier"
# sheet's summary-statistics table, carrying only lam (F) / kappa (G) values
# (no price/stdev/multiplier/volume data for this synthetic spot series).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simplified contract multiplier")
$ws.Activate()

$xlPasteFormats = -4122

# Row 15's formatting should match the existing table: column A uses the
# header/label style, F and G use the "Comma, 3 decimals" style already used
# throughout columns F/G. Copy formats down from the last populated row
# instead of hard-coding a style index.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial($xlPasteFormats)

$ws.Range("G14").Copy()
$ws.Range("F15").PasteSpecial($xlPasteFormats)
$ws.Range("G15").PasteSpecial($xlPasteFormats)

$ws.Range("A15").Value = "WTI-spot"
$ws.Range("F15").Value = 0.015
$ws.Range("G15").Value = 0.001

$excel.CutCopyMode = $false
$ws.Range("F16").Select()
